# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting the refreshed data snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 274
$ws.Range("F3").Value = 648
$ws.Range("F6").Value = 2839
$ws.Range("F8").Value = 59
$ws.Range("F9").Value = 30
$ws.Range("F10").Value = 585
$ws.Range("F14").Value = 6017
$ws.Range("F15").Value = 641
$ws.Range("F16").Value = 1051
$ws.Range("F17").Value = 16
$ws.Range("F18").Value = 244
$ws.Range("F19").Value = 175
$ws.Range("F21").Value = 558
$ws.Range("F22").Value = 6
$ws.Range("F23").Value = 45
$ws.Range("F25").Value = 132
$ws.Range("F26").Value = 1326
$ws.Range("F29").Value = 54
$ws.Range("F30").Value = 2073
$ws.Range("F32").Value = 10
$ws.Range("F33").Value = 360
$ws.Range("F35").Value = 3324

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 94
$ws.Range("F22").Value = 356
$ws.Range("F28").Value = 153
$ws.Range("F29").Value = 217
$ws.Range("F32").Value = 197
$ws.Range("F35").Value = 43

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2602
$ws.Range("F6").Value = 1146
$ws.Range("F8").Value = 1500
$ws.Range("F10").Value = 118
$ws.Range("F12").Value = 663

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2602
$ws.Range("F5").Value = 1146
$ws.Range("F6").Value = 1500
$ws.Range("F8").Value = 118
$ws.Range("F9").Value = 274
$ws.Range("F10").Value = 648
$ws.Range("F11").Value = 2839
$ws.Range("F12").Value = 59
$ws.Range("F13").Value = 663
$ws.Range("F14").Value = 585
$ws.Range("F15").Value = 94
$ws.Range("F19").Value = 6017
$ws.Range("F21").Value = 641
$ws.Range("F22").Value = 1051
$ws.Range("F23").Value = 16
$ws.Range("F24").Value = 244
$ws.Range("F25").Value = 175
$ws.Range("F27").Value = 558
$ws.Range("F33").Value = 356
$ws.Range("F37").Value = 153
$ws.Range("F38").Value = 217
$ws.Range("F40").Value = 54
$ws.Range("F43").Value = 2073
$ws.Range("F44").Value = 197
$ws.Range("F45").Value = 43
$ws.Range("F47").Value = 360
$ws.Range("F49").Value = 3324

Write-Output "Updated F-column values across 展览, 演出, 本地生活, 全部类型 sheets."
